# Regenerate the "K" column (G) of save_data with new simulated strikeout values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2 = 3; 3 = 1; 4 = 1; 5 = 2; 6 = 2; 7 = 1; 8 = 0; 9 = 1; 10 = 3;
    11 = 1; 12 = 1; 13 = 0; 14 = 0; 15 = 0; 16 = 0; 17 = 1; 18 = 2; 19 = 3;
    20 = 3; 21 = 1; 22 = 1; 23 = 2; 24 = 1; 25 = 3; 26 = 1; 27 = 0; 28 = 0;
    29 = 0; 30 = 2; 31 = 2; 32 = 1; 33 = 0; 34 = 2; 35 = 2; 36 = 2; 37 = 2;
    38 = 0; 39 = 2; 40 = 1; 41 = 0; 42 = 1; 43 = 0; 44 = 1; 45 = 2; 46 = 2;
    47 = 0; 48 = 2; 49 = 1; 50 = 0; 51 = 1; 52 = 0; 53 = 2; 54 = 1; 55 = 0;
    56 = 0; 57 = 1; 58 = 2; 59 = 1; 60 = 1; 61 = 0; 62 = 1; 63 = 0; 64 = 1;
    65 = 0; 66 = 1; 67 = 1; 68 = 1; 69 = 2; 70 = 3; 71 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
